$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.405.77"
$ws.Range("D3").Value = "1.937.93"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7477"
$ws.Range("E5").Value = "  +4.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "245.17"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.60"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3169"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06979"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7802"
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07995"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "1.938.87"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.369"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.49"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.41"
$ws.Range("E16").Value = "  -4.13%  "
$ws.Range("D17").Value = "30.397.05"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.08"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007922"
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.746"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "2.193.07"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.680"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.511"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.72"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.99"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1329"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.243"
$ws.Range("E29").Value = "  -4.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.358"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.511"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.362"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05153"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7472"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.787"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01949"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.802"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "77.94"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.415"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4454"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.963"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8343"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.37"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.749"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.446"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "985.17"
$ws.Range("E49").Value = "  +11.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.21"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06024"
$ws.Range("E51").Value = "  -0.53%  "

Write-Output "Updated cryptos list"
